$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps its text formatting so values like
# "25.935.89" or "0.9994" are not coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '25.935.89'
$ws.Range("E2").Value = '  -0.71%  '
$ws.Range("D3").Value = '1.745.26'
$ws.Range("E3").Value = '  -0.18%  '
$ws.Range("D4").Value = '0.9994'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '247.45'
$ws.Range("E5").Value = '  +4.63%  '
$ws.Range("D6").Value = '0.9995'
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").Value = '0.5043'
$ws.Range("D8").Value = '0.2741'
$ws.Range("E8").Value = '  -2.48%  '
$ws.Range("D9").Value = '0.06187'
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("D10").Value = '1.748.42'
$ws.Range("E10").Value = '  -0.02%  '
$ws.Range("D11").Value = '0.07258'
$ws.Range("E11").Value = '  +1.14%  '
$ws.Range("D12").Value = '0.6546'
$ws.Range("E12").Value = '  +1.46%  '
$ws.Range("E13").Value = '  -1.40%  '
$ws.Range("D14").Value = '4.641'
$ws.Range("E14").Value = '  +0.21%  '
$ws.Range("E15").Value = '  -1.14%  '
$ws.Range("D16").Value = '0.9998'
$ws.Range("D17").Value = '0.9990'
$ws.Range("E17").Value = '  -0.08%  '
$ws.Range("D18").Value = '25.962.88'
$ws.Range("E18").Value = '  -0.26%  '
$ws.Range("D19").Value = '11.84'
$ws.Range("E19").Value = '  +0.73%  '
$ws.Range("D20").Value = '0.000006824'
$ws.Range("E20").Value = '  +1.29%  '
$ws.Range("D21").Value = '1.969.96'
$ws.Range("D22").Value = '4.362'
$ws.Range("E22").Value = '  +0.93%  '
$ws.Range("D23").Value = '8.692'
$ws.Range("E23").Value = '  -0.34%  '
$ws.Range("D24").Value = '5.400'
$ws.Range("E24").Value = '  +3.33%  '
$ws.Range("D25").Value = '136.73'
$ws.Range("E25").Value = '  -2.08%  '
$ws.Range("E26").Value = '  -1.21%  '
$ws.Range("D27").Value = '15.25'
$ws.Range("E27").Value = '  -0.26%  '
$ws.Range("D28").Value = '1.778'
$ws.Range("E28").Value = '  -1.77%  '
$ws.Range("D29").Value = '105.44'
$ws.Range("E29").Value = '  +0.55%  '
$ws.Range("D30").Value = '3.904'
$ws.Range("E30").Value = '  +2.73%  '
$ws.Range("D31").Value = '0.08230'
$ws.Range("D32").Value = '3.639'
$ws.Range("D33").Value = '0.04676'
$ws.Range("E33").Value = '  +1.00%  '
$ws.Range("D34").Value = '2.653'
$ws.Range("E34").Value = '  +0.17%  '
$ws.Range("D35").Value = '0.9937'
$ws.Range("E35").Value = '  -1.52%  '
$ws.Range("D36").Value = '0.6191'
$ws.Range("E36").Value = '  -2.27%  '
$ws.Range("D37").Value = '2.757'
$ws.Range("E37").Value = '  +2.02%  '
$ws.Range("D38").Value = '0.01611'
$ws.Range("E38").Value = '  -0.73%  '
$ws.Range("D39").Value = '1.931'
$ws.Range("E39").Value = '  -2.24%  '
$ws.Range("D40").Value = '0.9994'
$ws.Range("E40").Value = '  -0.02%  '
$ws.Range("D41").Value = '100.00'
$ws.Range("E41").Value = '  -2.40%  '
$ws.Range("D42").Value = '0.3923'
$ws.Range("E42").Value = '  -0.06%  '
$ws.Range("D43").Value = '0.7583'
$ws.Range("E43").Value = '  +0.68%  '
$ws.Range("D44").Value = '5.008'
$ws.Range("E44").Value = '  -0.93%  '
$ws.Range("E45").Value = '  -0.85%  '
$ws.Range("D46").Value = '6.298'
$ws.Range("E46").Value = '  -0.91%  '
$ws.Range("D47").Value = '55.51'
$ws.Range("E47").Value = '  +1.71%  '
$ws.Range("D48").Value = '0.05260'
$ws.Range("E48").Value = '  -1.67%  '
$ws.Range("D49").Value = '30.63'
$ws.Range("E49").Value = '  -1.29%  '
$ws.Range("B50").Value = 'Decentraland'
$ws.Range("C50").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D50").Value = '0.3433'
$ws.Range("E50").Value = '  -1.39%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '7.521'
$ws.Range("E51").Value = '  -0.61%  '

Write-Host "Updated cryptos list"